$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily cryptos list refresh (GitHub Actions bot): update Price (D) and
# Volume(1h) (E) columns for every coin row, and fix the RenderToken /
# VeChain rows (40-41), which swapped rank order, including their
# Coin name (B) and Link (C) cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.800.90"
$ws.Range("E2").Value = "  +0.45%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.749.92"
$ws.Range("E3").Value = "  +0.37%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.38"
$ws.Range("E5").Value = "  -0.20%  "

$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5078"
$ws.Range("E7").Value = "  +3.45%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.64"
$ws.Range("E8").Value = "  -2.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2694"
$ws.Range("E9").Value = "  +7.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06194"
$ws.Range("E10").Value = "  +4.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.750.26"
$ws.Range("E11").Value = "  +0.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06927"
$ws.Range("E12").Value = "  +2.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.52"
$ws.Range("E13").Value = "  +6.47%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6254"
$ws.Range("E14").Value = "  +9.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.480"
$ws.Range("E15").Value = "  +0.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.91"
$ws.Range("E16").Value = "  +1.22%  "

$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("E18").Value = "  -0.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "25.813.32"
$ws.Range("E19").Value = "  +0.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.68"
$ws.Range("E20").Value = "  +1.82%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006696"
$ws.Range("E21").Value = "  +2.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.977.08"
$ws.Range("E22").Value = "  +0.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.058"
$ws.Range("E23").Value = "  +2.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.247"
$ws.Range("E24").Value = "  +4.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.158"
$ws.Range("E25").Value = "  +2.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "136.68"
$ws.Range("E26").Value = "  +0.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.17"
$ws.Range("E27").Value = "  +4.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.455"
$ws.Range("E28").Value = "  -1.92%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.764"
$ws.Range("E29").Value = "  -3.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.49"
$ws.Range("E30").Value = "  +1.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08252"
$ws.Range("E31").Value = "  +2.92%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.709"
$ws.Range("E32").Value = "  -1.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.412"
$ws.Range("E33").Value = "  +2.94%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04435"
$ws.Range("E34").Value = "  +0.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.000"
$ws.Range("E35").Value = "  -0.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.649"
$ws.Range("E36").Value = "  +1.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9997"
$ws.Range("E37").Value = "  +0.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6031"
$ws.Range("E38").Value = "  +0.66%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.693"
$ws.Range("E39").Value = "  -0.25%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01562"
$ws.Range("E40").Value = "  +5.24%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.953"
$ws.Range("E41").Value = "  -3.94%  "

$ws.Range("E42").Value = "  -0.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.28"
$ws.Range("E43").Value = "  -2.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3849"
$ws.Range("E44").Value = "  +3.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7504"
$ws.Range("E45").Value = "  -1.75%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.897"
$ws.Range("E46").Value = "  -5.60%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05501"
$ws.Range("E47").Value = "  +7.59%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1101"
$ws.Range("E48").Value = "  +2.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.973"
$ws.Range("E49").Value = "  +1.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.17"
$ws.Range("E50").Value = "  +0.26%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.84"
$ws.Range("E51").Value = "  +0.58%  "
